$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "201VP00023"
$ws.Range("A3").Value = "201VP00022"
$ws.Range("A4").Value = "201VP00026"

$ws.Range("E8").Select()
